$wb = $excel.ActiveWorkbook

# --- Sheet4 rename: "Transposed x 1000 + IEA report" -> "... IEA country" ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Transposed x 1000 + IEA country"

# --- Sheet1 view: scroll to A65, select whole row 196 ---
$ws1.Activate()
$ws1.Rows.Item(196).Select()
$excel.ActiveWindow.ScrollRow = 65
$excel.ActiveWindow.ScrollColumn = 1

# --- Sheet2 view / Sheet3 view: left unchanged (already A2:M229 / B1) ---

# --- Sheet4 view: freeze panes at B2 (row1/col1 frozen), selection L31 in bottomRight pane ---
$ws4.Activate()
$ws4.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws4.Range("L31").Select()
